$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 ---
# B2 "1234" -> "5634" (still a text value, so force text format before
# writing the numeric-looking string, otherwise Excel auto-converts it
# to a number).
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "5634"

$ws.Range("C2").Value = "Ninja ZX25R"
$ws.Range("D2").Value = "45.000.000"
# E2 stays "34" - unchanged

# --- Add new row 3 ---
# A3 mirrors A2's style (s="1") and numeric value 1.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 1

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2468"

$ws.Range("C3").Value = "Benelli TRK 502"
$ws.Range("D3").Value = "54.000.000"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "40"
